# Applies the "Updated symbol list" edit: refreshed prices and a one-row
# upward rotation of the exchange-token listings in rows 9-26.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("D6").Value = "'3.449"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.8092"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8874"
$ws.Range("D8").Style = "Normal"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1446"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07366"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03022"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03072"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09405"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'3.936"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13MCDexMCBBestin24h"
$ws.Range("D15").Value = "'0.001581"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.04811"
$ws.Range("D16").Style = "Normal"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005849"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16OneONE"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006151"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.005104"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.0009977"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.0001500"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.749"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "KuCoinToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D23").Value = "'6.298"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22KuCoinTokenKCS"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.191"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").Value = "'0.3276"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").Value = "'0.1320"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("D27").Value = "'0.0003001"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03902"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006791"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1068"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002520"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007377"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005640"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.3800"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.1732"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
